$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Min Cost Agent")
$ws.Name = "Max Profit Agent"
